# Add a description for metadata in cell E3, and move the active selection to E4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Data about location"

$ws.Range("E4").Select()
